# Add variant of cancer status derived variable (Ca07a / cancer_status_v2)
# Inserts a new row just above the existing "Ca08" row (current row 25),
# shifting the remainder of the table down by one row, then fills in the
# new row's values and refreshes the table range / selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 25 - this shifts rows 25..130 down to 26..131
# and updates the sheet dimension automatically.
$ws.Rows.Item(25).Insert()

# Grow Table1 (was A1:E130) to cover the newly inserted row (A1:E131).
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E131"))

# Populate the new row with the "Ca07a" / cancer_status_v2 variable.
$ws.Range("A25").Value = "Ca07a"
$ws.Range("B25").Value = "cancer_status_v2"
$ws.Range("C25").Value = "Cancer"
$ws.Range("D25").Value = "Same as Ca07 except do not combine stable and responding"

# Move the view/selection to the newly added cell, matching the edit.
$ws.Application.ActiveWindow.ScrollRow = 17
$ws.Application.ActiveWindow.TopLeftCell = $ws.Range("A17")
[void]$ws.Range("D25").Select()
